$d = $word.ActiveDocument

# Original contact line (single run) and its replacement pieces.
$oldText  = "Williamsport, PA | kolakowski.davis@gmail.com | www.linkedin.com/in/davis-kolakowski | (973) 974-0088"
$prefix   = "Williamsport, PA | "
$newEmail = "contact@daviskolakowski.com"
$suffix   = " | www.linkedin.com/in/davis-kolakowski"
$newText  = $prefix + $newEmail + $suffix

# Locate the old text anywhere in the document body.
$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $baseStart = $rng.Start

    # Replace the whole phrase in-place; since this happens inside the
    # single original run, the run keeps its existing character formatting.
    $rng.Text = $newText

    $p1End = $baseStart + $prefix.Length
    $p2End = $p1End + $newEmail.Length
    $p3End = $baseStart + $newText.Length

    # Touch (and restore) formatting on the email + linkedin segments so
    # Word splits them into their own runs, each carrying the same
    # character formatting as the original run.
    $r2 = $d.Range($p1End, $p2End)
    $r2.Bold = $true
    $r2.Bold = $false

    $r3 = $d.Range($p2End, $p3End)
    $r3.Bold = $true
    $r3.Bold = $false
}
